# Update "Bestand" (column D) values on the "Lagerbestand M0129" sheet
# as part of translating notifications and permission error handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = 2575
    9   = 2491
    14  = 4768
    22  = 2017
    29  = 6017
    30  = -235
    33  = 29
    40  = -20
    41  = 5473
    42  = 3771
    43  = 232
    45  = 3100
    46  = 151
    48  = 72
    51  = -400
    52  = 244
    53  = 196
    56  = 1786
    63  = 47
    75  = 211
    76  = 4594
    82  = 216
    90  = 686
    91  = 691
    92  = 857
    93  = 629
    97  = 868
    99  = 463
    118 = 2
    124 = 308
    125 = 62
    127 = 403
    155 = 2242
    160 = 702
    181 = 38
    182 = -40
    207 = -2008.5
    225 = 42
    230 = -500
    231 = 20
    238 = -2000
    242 = -10
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 4).Value = $updates[$row]
}
